$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the Arabic message text in D4 ("...one of its two fields is empty"
#    -> "...the original-text field is empty"), replacing the shared string
#    in place so every other cell referencing it is unaffected.
$null = $ws.Cells.Replace("سيتم تجاهل الصف الذي يكون أحد حقليه فارغاً.", "سيتم تجاهل الصف الذي يكون فيه حقل النص الأصلي فارغاً.")

# 2. Row 4 height shrinks from 15 to 13.8 (matching the other data rows).
$ws.Rows.Item(4).RowHeight = 13.8

# 3. Update the view: scroll so column B is the left-most visible column,
#    and move the active selection to the single cell D7.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$null = $ws.Range("D7").Select()
